$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 352, shifting existing rows 352:369 down to 353:370
$ws.Rows.Item(352).Insert()

# Fill the new row 352 with the new record (copy of old row 352's static fields,
# with updated Fecha/Precio minimo/Precio maximo/Precio promedio ponderado/Precio $/Kg)
$ws.Cells.Item(352, 1).Value = 3
$ws.Cells.Item(352, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(352, 3).Value = "Coquimbo"
$ws.Cells.Item(352, 4).Value = 44753
$ws.Cells.Item(352, 5).Value = 5
$ws.Cells.Item(352, 6).Value = 100112009
$ws.Cells.Item(352, 7).Value = "Acelga"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 230
$ws.Cells.Item(352, 11).Value = 3300
$ws.Cells.Item(352, 12).Value = 3500
$ws.Cells.Item(352, 13).Value = 3404
$ws.Cells.Item(352, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(352, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(352, 16).Value = 567
$ws.Cells.Item(352, 17).Value = 6
$ws.Cells.Item(352, 18).Value = "Hortaliza"
